$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column).
# Excel's native "insert column" shifts existing columns/cells to the right
# and inherits formatting from the column to the left.
$ws.Columns("N").Insert() | Out-Null

# The newly inserted column inherits column M's raw width value when Excel
# performs a real column-insert; approximate that width (closest value the
# engine can store) and mark it as an explicit (non bestFit) custom width.
$ws.Columns("N").ColumnWidth = 9.9

# Make "Repayment schedule" the active sheet/tab and select cell R6,
# matching the new selection saved in the workbook.
$ws.Activate() | Out-Null
$ws.Range("R6").Select() | Out-Null
